# Updated RAD EL-Motor Fuel Tax.
# The Katalon RAD regression test was re-run, so the "Date" column
# (column B) timestamps for the three recorded test steps need to be
# refreshed to the new execution times.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Mon Apr 08 18:07:58 EDT 2024"
$ws.Range("B3").Value = "Mon Apr 08 18:08:10 EDT 2024"
$ws.Range("B4").Value = "Mon Apr 08 18:08:23 EDT 2024"
